$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New regenerated s_val data (filtered save games), rows 2-10, columns B:G
$data = @{
    2  = @(3.272327238179451, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 5.582307763322248)
    3  = @(3.272327238179451, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 6.15379541431027)
    4  = @(1.445647641019636, 0.3048912486333797, 0.7210945179870265, 0.5333859586016987, 3.005019366241741)
    5  = @(0.6545652718822623, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 2.964545797025059)
    6  = @(3.272327238179451, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 6.15379541431027)
    7  = @(3.272327238179451, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 6.15379541431027)
    8  = @(0.1169995834814548, 0.3048912486333797, 0.1496068669990043, 0.5333859586016987, 1.104883657715537)
    9  = @(0.1169995834814548, 0.04103571897497393, 0.7210945179870265, 0.5333859586016987, 1.412515779045154)
    10 = @(0.2881169905109251, 0.3048912486333797, 3.223369029078222, 0.5333859586016987, 4.349763226824225)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]  # B
    $ws.Cells.Item($row, 3).Value = $vals[1]  # C
    $ws.Cells.Item($row, 4).Value = $vals[2]  # D
    $ws.Cells.Item($row, 5).Value = $vals[3]  # E
    $ws.Cells.Item($row, 7).Value = $vals[4]  # G
}
